# The author removed the "Beslutstabeller/beslutsträd" section (4 slides:
# "Beslutstabeller/beslutsträd", "Beslutstabell", "Beslutsträd", "Testfall")
# which previously sat at slide-show positions 23-26. Deleting them shifts
# every following slide up by four positions, which also explains the
# renumbered relationship ids (notesMaster / custDataLst tags) seen in the
# diff. Deleting slides (and their notes/tags) is handled automatically by
# PowerPoint when a Slide is removed.

$p = $ppt.ActivePresentation

# Delete from the highest index down so the remaining indices never shift
# out from under us.
$p.Slides.Item(26).Delete()   # "Testfall"
$p.Slides.Item(25).Delete()   # "Beslutsträd"
$p.Slides.Item(24).Delete()   # "Beslutstabell"
$p.Slides.Item(23).Delete()   # "Beslutstabeller/beslutsträd"
